$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.854.18'
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '2.326.86'
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.02'
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.55'
$ws.Range("E6").Value = '  -2.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.508'
$ws.Range("E7").Value = '  -4.71%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.505'
$ws.Range("E9").Value = '  -4.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.44'
$ws.Range("E10").Value = '  -5.87%  '
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("E12").Value = '  -2.26%  '
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.75'
$ws.Range("E14").Value = '  -3.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.71'
$ws.Range("E15").Value = '  +3.96%  '
$ws.Range("D16").Value = '2.337.77'
$ws.Range("E16").Value = '  +1.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.822'
$ws.Range("E17").Value = '  +1.54%  '
$ws.Range("D18").Value = '42.779.73'
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0905'
$ws.Range("E19").Value = '  -2.35%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.14'
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.57'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.09'
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.43'
$ws.Range("E23").Value = '  -3.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("E24").Value = '  -2.01%  '
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  +2.23%  '
$ws.Range("E28").Value = '  -0.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  -2.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.71'
$ws.Range("E30").Value = '  -6.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.20'
$ws.Range("E31").Value = '  -4.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '161.03'
$ws.Range("E32").Value = '  -4.32%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.05'
$ws.Range("E34").Value = '  -4.41%  '
$ws.Range("E35").Value = '  +2.84%  '
$ws.Range("E36").Value = '  -3.35%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0718'
$ws.Range("E37").Value = '  -3.50%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.11'
$ws.Range("E38").Value = '  -5.67%  '
$ws.Range("E39").Value = '  -5.61%  '
$ws.Range("E40").Value = '  -2.71%  '
$ws.Range("E41").Value = '  -4.69%  '
$ws.Range("E42").Value = '  -3.30%  '
$ws.Range("E43").Value = '  -7.12%  '
$ws.Range("D44").Value = '2.009.23'
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("E45").Value = '  -4.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.49'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.15'
$ws.Range("E47").Value = '  +1.57%  '
$ws.Range("E48").Value = '  -4.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.44'
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("D51").Value = '2.555.89'
$ws.Range("E51").Value = '  +0.94%  '
